$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'30.552.21"
$ws.Range("E2").Value = "  -0.54%  "

# Row 3
$ws.Range("D3").Value = "'1.881.48"
$ws.Range("E3").Value = "  -0.38%  "

# Row 4
$ws.Range("D4").Value = "'0.9993"
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").Value = "'246.37"
$ws.Range("E5").Value = "  -0.68%  "

# Row 6
$ws.Range("D6").Value = "'0.9996"
$ws.Range("E6").Value = "  +0.00%  "

# Row 7
$ws.Range("D7").Value = "'0.4724"
$ws.Range("E7").Value = "  -0.36%  "

# Row 8
$ws.Range("D8").Value = "'0.2888"
$ws.Range("E8").Value = "  -1.42%  "

# Row 9
$ws.Range("D9").Value = "'0.06538"
$ws.Range("E9").Value = "  +0.09%  "

# Row 10
$ws.Range("D10").Value = "'22.05"
$ws.Range("E10").Value = "  +0.25%  "

# Row 11
$ws.Range("D11").Value = "'101.08"
$ws.Range("E11").Value = "  +4.18%  "

# Row 12
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").Value = "'0.7591"
$ws.Range("E12").Value = "  +3.14%  "

# Row 13
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.07822"
$ws.Range("E13").Value = "  +0.30%  "

# Row 14
$ws.Range("D14").Value = "'1.877.80"
$ws.Range("E14").Value = "  -0.58%  "

# Row 15
$ws.Range("D15").Value = "'5.248"
$ws.Range("E15").Value = "  +0.00%  "

# Row 16
$ws.Range("D16").Value = "'285.12"
$ws.Range("E16").Value = "  +0.19%  "

# Row 17
$ws.Range("D17").Value = "'30.531.52"
$ws.Range("E17").Value = "  -0.56%  "

# Row 19
$ws.Range("D19").Value = "'0.000007529"
$ws.Range("E19").Value = "  -0.05%  "

# Row 20
$ws.Range("D20").Value = "'0.9996"

# Row 21
$ws.Range("D21").Value = "'2.120.09"
$ws.Range("E21").Value = "  -0.87%  "

# Row 22
$ws.Range("D22").Value = "'5.374"
$ws.Range("E22").Value = "  +0.73%  "

# Row 23
$ws.Range("D23").Value = "'0.9990"
$ws.Range("E23").Value = "  -0.02%  "

# Row 24
$ws.Range("D24").Value = "'6.393"
$ws.Range("E24").Value = "  +2.19%  "

# Row 25
$ws.Range("D25").Value = "'9.137"
$ws.Range("E25").Value = "  -1.08%  "

# Row 26
$ws.Range("D26").Value = "'163.09"
$ws.Range("E26").Value = "  -0.89%  "

# Row 27
$ws.Range("D27").Value = "'19.08"
$ws.Range("E27").Value = "  +0.79%  "

# Row 28
$ws.Range("D28").Value = "'1.915"
$ws.Range("E28").Value = "  -0.46%  "

# Row 29
$ws.Range("D29").Value = "'0.09700"
$ws.Range("E29").Value = "  -0.44%  "

# Row 30
$ws.Range("E30").Value = "  -1.23%  "

# Row 31
$ws.Range("D31").Value = "'1.495"
$ws.Range("E31").Value = "  +0.04%  "

# Row 32
$ws.Range("D32").Value = "'4.265"
$ws.Range("E32").Value = "  -0.95%  "

# Row 33
$ws.Range("D33").Value = "'4.192"
$ws.Range("E33").Value = "  +0.14%  "

# Row 34
$ws.Range("D34").Value = "'0.04839"
$ws.Range("E34").Value = "  -0.53%  "

# Row 35
$ws.Range("D35").Value = "'1.130"
$ws.Range("E35").Value = "  +0.25%  "

# Row 36
$ws.Range("D36").Value = "'0.6954"
$ws.Range("E36").Value = "  -0.21%  "

# Row 37
$ws.Range("D37").Value = "'2.772"
$ws.Range("E37").Value = "  +1.84%  "

# Row 38
$ws.Range("E38").Value = "  +0.29%  "

# Row 39
$ws.Range("D39").Value = "'2.865"
$ws.Range("E39").Value = "  +2.22%  "

# Row 40
$ws.Range("D40").Value = "'76.56"
$ws.Range("E40").Value = "  +0.72%  "

# Row 41
$ws.Range("D41").Value = "'6.358"
$ws.Range("E41").Value = "  -0.71%  "

# Row 42
$ws.Range("D42").Value = "'1.979"
$ws.Range("E42").Value = "  -1.56%  "

# Row 43
$ws.Range("D43").Value = "'0.4255"
$ws.Range("E43").Value = "  -0.16%  "

# Row 44
$ws.Range("D44").Value = "'0.9993"
$ws.Range("E44").Value = "  -0.03%  "

# Row 45
$ws.Range("D45").Value = "'0.8294"
$ws.Range("E45").Value = "  -0.74%  "

# Row 46
$ws.Range("D46").Value = "'101.34"
$ws.Range("E46").Value = "  -0.42%  "

# Row 47
$ws.Range("D47").Value = "'9.818"
$ws.Range("E47").Value = "  +3.08%  "

# Row 48
$ws.Range("D48").Value = "'7.055"
$ws.Range("E48").Value = "  +0.39%  "

# Row 49
$ws.Range("D49").Value = "'35.12"

# Row 50
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "'890.83"
$ws.Range("E50").Value = "  -3.13%  "

# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.05766"
$ws.Range("E51").Value = "  +0.27%  "
